$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Source data: the "defined range" this macro scans. Each tuple is
# (fund, NAV per, BM, BM tol). For every fund we work out the BM Bps
# Var (BM - NAV per) and flag it when the variance breaks tolerance.
# ------------------------------------------------------------------
$sourceRange = @(
    @{ Fund = 3123; NavPer = 100; BM = 130; Tol = 20 },
    @{ Fund = 3123; NavPer = 100; BM = 130; Tol = 20 },
    @{ Fund = 3123; NavPer = 100; BM = 130; Tol = 20 },
    @{ Fund = 3456; NavPer = 120; BM = 110; Tol = 40 },
    @{ Fund = 3456; NavPer = 120; BM = 110; Tol = 40 },
    @{ Fund = 111;  NavPer = -60; BM = -10; Tol = 30 },
    @{ Fund = 2323; NavPer = -60; BM = -10; Tol = 30 },
    @{ Fund = 3969; NavPer = -70; BM = -60; Tol = 30 }
)

# Headers for the dictionary we're about to build / print out.
$headers = @("fund", "NAV per", "BM", "BM tol", "BM Bps Var")
for ($h = 0; $h -lt $headers.Length; $h++) {
    $ws.Cells.Item(1, 11 + $h).Value2 = $headers[$h]
}

# Walk the defined range, build the fund -> variant dictionary, and
# drop each entry into the next free row under the headers.
$dictionary = @{}
$destRow = 2
foreach ($entry in $sourceRange) {
    $fund = $entry.Fund
    $navPer = $entry.NavPer
    $bm = $entry.BM
    $tol = $entry.Tol

    # Add fund + its variant to the dictionary if not already captured.
    if (-not $dictionary.ContainsKey($destRow)) {
        $dictionary[$destRow] = $fund
    }

    $ws.Cells.Item($destRow, 11).Value2 = $fund
    $ws.Cells.Item($destRow, 12).Value2 = $navPer
    $ws.Cells.Item($destRow, 13).Value2 = $bm
    $ws.Cells.Item($destRow, 14).Value2 = $tol
    $destRow = $destRow + 1
}

# BM Bps Var = BM - NAV per, written as one shared formula across the
# whole block so Excel stores it the same way AutoFill would.
$lastRow = $destRow - 1
$ws.Range("O2:O" + $lastRow).Formula = "=M2-L2"

# Highlight in yellow any fund whose variance breaks its tolerance.
for ($r = 2; $r -le $lastRow; $r++) {
    $variance = $ws.Cells.Item($r, 13).Value2 - $ws.Cells.Item($r, 12).Value2
    $tolerance = $ws.Cells.Item($r, 14).Value2
    if ([Math]::Abs($variance) -gt $tolerance) {
        $ws.Cells.Item($r, 15).Interior.Color = 65535
    }
}

$null = $ws.Range("L13").Select()
